$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Price column (D) values that look numeric need to be forced to text
# so Excel stores them verbatim (matching the source data format) instead of
# auto-converting to a Number and losing formatting such as trailing zeros.
$priceCells = @("D4", "D5", "D6", "D8", "D9", "D12", "D13", "D14", "D15", "D16", "D17", "D19", "D21", "D22", "D23", "D24", "D26", "D27", "D28", "D29", "D30", "D31", "D33", "D35", "D37", "D40", "D41", "D43", "D44", "D45", "D46", "D47", "D48", "D49")
foreach ($addr in $priceCells) {
    $ws.Range($addr).NumberFormat = "@"
}

$ws.Range("D4").Value = "1.0000"
$ws.Range("D5").Value = "240.10"
$ws.Range("D6").Value = "0.6288"
$ws.Range("D8").Value = "0.07652"
$ws.Range("D9").Value = "0.2915"
$ws.Range("D12").Value = "0.07735"
$ws.Range("D13").Value = "5.032"
$ws.Range("D14").Value = "0.6809"
$ws.Range("D15").Value = "0.00001057"
$ws.Range("D16").Value = "83.50"
$ws.Range("D17").Value = "6.199"
$ws.Range("D19").Value = "229.17"
$ws.Range("D21").Value = "1.000"
$ws.Range("D22").Value = "7.461"
$ws.Range("D23").Value = "1.001"
$ws.Range("D24").Value = "157.27"
$ws.Range("D26").Value = "8.424"
$ws.Range("D27").Value = "17.72"
$ws.Range("D28").Value = "1.382"
$ws.Range("D29").Value = "1.465"
$ws.Range("D30").Value = "0.05612"
$ws.Range("D31").Value = "4.130"
$ws.Range("D33").Value = "1.844"
$ws.Range("D35").Value = "0.6999"
$ws.Range("D37").Value = "0.01804"
$ws.Range("D40").Value = "6.468"
$ws.Range("D41").Value = "0.9073"
$ws.Range("D43").Value = "101.97"
$ws.Range("D44").Value = "66.06"
$ws.Range("D45").Value = "7.206"
$ws.Range("D46").Value = "0.00000000119"
$ws.Range("D47").Value = "0.4028"
$ws.Range("D48").Value = "0.1153"
$ws.Range("D49").Value = "8.996"

foreach ($addr in $priceCells) {
    $ws.Range($addr).Style = "Normal"
}

# --- Remaining cells (Volume % column, and Price values that already read as text)
$ws.Range("D2").Value = "29.438.37"
$ws.Range("E2").Value = "  +0.15%  "
$ws.Range("D3").Value = "1.850.12"
$ws.Range("E3").Value = "  +0.10%  "
$ws.Range("E4").Value = "  +0.05%  "
$ws.Range("E5").Value = "  +0.00%  "
$ws.Range("E6").Value = "  -0.15%  "
$ws.Range("E7").Value = "  +0.02%  "
$ws.Range("E8").Value = "  +0.23%  "
$ws.Range("E9").Value = "  -0.51%  "
$ws.Range("E10").Value = "  +1.01%  "
$ws.Range("D11").Value = "2.149.02"
$ws.Range("E11").Value = "  +15.92%  "
$ws.Range("E12").Value = "  -0.06%  "
$ws.Range("E13").Value = "  +0.56%  "
$ws.Range("E15").Value = "  -5.64%  "
$ws.Range("E16").Value = "  -0.32%  "
$ws.Range("E17").Value = "  +0.35%  "
$ws.Range("D18").Value = "29.501.76"
$ws.Range("E18").Value = "  +0.31%  "
$ws.Range("E19").Value = "  +0.09%  "
$ws.Range("E21").Value = "  +0.05%  "
$ws.Range("E22").Value = "  -0.41%  "
$ws.Range("E23").Value = "  +0.01%  "
$ws.Range("E24").Value = "  -0.06%  "
$ws.Range("E25").Value = "  -0.82%  "
$ws.Range("E26").Value = "  +0.92%  "
$ws.Range("E27").Value = "  +0.52%  "
$ws.Range("E28").Value = "  +6.33%  "
$ws.Range("E29").Value = "  -0.17%  "
$ws.Range("E30").Value = "  +0.24%  "
$ws.Range("E31").Value = "  +0.34%  "
$ws.Range("E32").Value = "  +0.82%  "
$ws.Range("E33").Value = "  -0.37%  "
$ws.Range("E34").Value = "  +0.79%  "
$ws.Range("E35").Value = "  -1.58%  "
$ws.Range("E36").Value = "  +0.23%  "
$ws.Range("E37").Value = "  -0.05%  "
$ws.Range("D38").Value = "1.231.40"
$ws.Range("E38").Value = "  -0.72%  "
$ws.Range("E39").Value = "  -1.42%  "
$ws.Range("E40").Value = "  +1.04%  "
$ws.Range("E41").Value = "  +0.37%  "
$ws.Range("E44").Value = "  +0.05%  "
$ws.Range("E45").Value = "  +0.73%  "
$ws.Range("E46").Value = "  +1.13%  "
$ws.Range("E47").Value = "  +0.28%  "
$ws.Range("E48").Value = "  +2.98%  "
$ws.Range("E49").Value = "  -0.54%  "
$ws.Range("E50").Value = "  -0.35%  "
$ws.Range("E51").Value = "  +0.01%  "
